$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J6").Value = 3.4
$ws.Range("Q6").Value = 2.5
$ws.Range("R6").Value = 1.5
$ws.Range("U6").Value = 2.1
$ws.Range("V6").Value = 1.67
$ws.Range("W6").Value = 6.5
$ws.Range("Y6").Value = 11
$ws.Range("AE6").Value = 19
$ws.Range("AK6").Value = 29
$ws.Range("AU6").Value = 9
$ws.Range("AX6").Value = 19
$ws.Range("AY6").Value = 34
$ws.Range("BB6").Value = 301
$ws.Range("I7").Value = 3.7
$ws.Range("K7").Value = 1.83
$ws.Range("AA7").Value = 23
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("O8").Value = 1.29
$ws.Range("P8").Value = 3.5
$ws.Range("Q8").Value = 1.98
$ws.Range("R8").Value = 1.88
$ws.Range("G15").Value = 1.65
$ws.Range("H15").Value = 3.5
$ws.Range("I15").Value = 4.75
$ws.Range("J15").Value = 2.3
$ws.Range("K15").Value = 2.1
$ws.Range("L15").Value = 5.5
$ws.Range("Q15").Value = 2.08
$ws.Range("R15").Value = 1.73
$ws.Range("U15").Value = 2
$ws.Range("V15").Value = 1.73
$ws.Range("X15").Value = 7.5
$ws.Range("Z15").Value = 12
$ws.Range("AD15").Value = 7
$ws.Range("AE15").Value = 19
$ws.Range("AG15").Value = 12
$ws.Range("AH15").Value = 26
$ws.Range("AI15").Value = 17
$ws.Range("AJ15").Value = 51
$ws.Range("AK15").Value = 41
$ws.Range("AM15").Value = 1000
$ws.Range("AN15").Value = 3.6
$ws.Range("AO15").Value = 9
$ws.Range("AQ15").Value = 29
$ws.Range("AW15").Value = 6.5
$ws.Range("AX15").Value = 29
$ws.Range("AZ15").Value = 101
$ws.Range("BA15").Value = 126
$ws.Range("BB15").Value = 301
$ws.Range("G16").Value = 2.35
$ws.Range("I16").Value = 2.9
$ws.Range("J16").Value = 3.2
$ws.Range("K16").Value = 2.05
$ws.Range("O16").Value = 1.36
$ws.Range("P16").Value = 3
$ws.Range("Q16").Value = 2.15
$ws.Range("R16").Value = 1.67
$ws.Range("U16").Value = 1.83
$ws.Range("V16").Value = 1.83
$ws.Range("AE16").Value = 15
$ws.Range("AG16").Value = 8.5
$ws.Range("AH16").Value = 15
$ws.Range("AR16").Value = 67
$ws.Range("AV16").Value = 51
$ws.Range("G36").Value = 2.87
$ws.Range("H36").Value = 3.1
$ws.Range("J36").Value = 3.4
$ws.Range("K36").Value = 2.02
$ws.Range("L36").Value = 3
$ws.Range("N36").Value = 9.4
$ws.Range("O36").Value = 1.28
$ws.Range("P36").Value = 3.05
$ws.Range("Q36").Value = 1.87
$ws.Range("R36").Value = 1.83
$ws.Range("S36").Value = 1.4
$ws.Range("T36").Value = 2.52
$ws.Range("W36").Value = 9.75
$ws.Range("X36").Value = 16
$ws.Range("AA36").Value = 23
$ws.Range("AB36").Value = 29
$ws.Range("AC36").Value = 9.75
$ws.Range("AD36").Value = 6.1
$ws.Range("AG36").Value = 8.25
$ws.Range("AH36").Value = 12
$ws.Range("AJ36").Value = 25
$ws.Range("AK36").Value = 19.5
$ws.Range("AL36").Value = 27
$ws.Range("AN36").Value = 4.8
$ws.Range("AO36").Value = 15.5
$ws.Range("AP36").Value = 21
$ws.Range("AQ36").Value = 70
$ws.Range("AR36").Value = 100
$ws.Range("AS36").Value = 250
$ws.Range("AT36").Value = 2.5
$ws.Range("AV36").Value = 55
$ws.Range("AX36").Value = 13
$ws.Range("AY36").Value = 21
$ws.Range("AZ36").Value = 55
$ws.Range("BA36").Value = 90
$ws.Range("G37").Value = 2.22
$ws.Range("H37").Value = 3.2
$ws.Range("I37").Value = 3.05
$ws.Range("K37").Value = 2.05
$ws.Range("L37").Value = 3.6
$ws.Range("M37").Value = 1.02
$ws.Range("N37").Value = 7.1
$ws.Range("O37").Value = 1.37
$ws.Range("P37").Value = 2.65
$ws.Range("Q37").Value = 2.07
$ws.Range("R37").Value = 1.6
$ws.Range("S37").Value = 1.42
$ws.Range("T37").Value = 2.45
$ws.Range("U37").Value = 1.85
$ws.Range("V37").Value = 1.75
$ws.Range("W37").Value = 6.8
$ws.Range("X37").Value = 10
$ws.Range("Y37").Value = 9.25
$ws.Range("Z37").Value = 21
$ws.Range("AB37").Value = 35
$ws.Range("AC37").Value = 8.25
$ws.Range("AD37").Value = 6.2
$ws.Range("AE37").Value = 16
$ws.Range("AF37").Value = 90
$ws.Range("AG37").Value = 8.25
$ws.Range("AH37").Value = 15
$ws.Range("AI37").Value = 11.25
$ws.Range("AJ37").Value = 37
$ws.Range("AK37").Value = 29
$ws.Range("AL37").Value = 40
$ws.Range("AM37").Value = 800
$ws.Range("AN37").Value = 4
$ws.Range("AO37").Value = 11.5
$ws.Range("AP37").Value = 22
$ws.Range("AR37").Value = 90
$ws.Range("AS37").Value = 300
$ws.Range("AT37").Value = 2.42
$ws.Range("AU37").Value = 7.4
$ws.Range("AV37").Value = 75
$ws.Range("AW37").Value = 4.8
$ws.Range("AX37").Value = 17
$ws.Range("AY37").Value = 26
$ws.Range("AZ37").Value = 80
$ws.Range("BA37").Value = 120
